$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Copies" header in F1, matching the style used by the other headers
# but centered without border (new style).
$ws.Range("F1").Value = "Copies"
$ws.Range("F1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("F1").Font.Bold = $true

# Fill in the copies count column
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 7
$ws.Range("F6").Value = 2
